$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue "D2" "29.908.70"
$ws.Range("E2").Value = "  -0.92%  "
Set-TextValue "D3" "1.897.14"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  -0.11%  "
Set-TextValue "D5" "0.7569"
$ws.Range("E5").Value = "  +2.25%  "
Set-TextValue "D6" "239.90"
$ws.Range("E6").Value = "  -1.74%  "
Set-TextValue "D7" "0.9998"
$ws.Range("E7").Value = "  -0.10%  "
Set-TextValue "D8" "0.3047"
$ws.Range("E8").Value = "  -2.73%  "
Set-TextValue "D9" "25.51"
$ws.Range("E9").Value = "  -5.36%  "
Set-TextValue "D10" "0.06846"
$ws.Range("E10").Value = "  -1.47%  "
Set-TextValue "D11" "0.07977"
$ws.Range("E11").Value = "  -0.16%  "
Set-TextValue "D12" "0.7464"
$ws.Range("E12").Value = "  -4.10%  "
Set-TextValue "D13" "1.891.93"
$ws.Range("E13").Value = "  -1.10%  "
Set-TextValue "D14" "5.189"
$ws.Range("E14").Value = "  -1.58%  "
Set-TextValue "D15" "91.28"
$ws.Range("E15").Value = "  -0.38%  "
Set-TextValue "D16" "29.903.03"
$ws.Range("E16").Value = "  -0.86%  "
Set-TextValue "D17" "13.93"
$ws.Range("E17").Value = "  -2.15%  "
Set-TextValue "D18" "5.954"
$ws.Range("E18").Value = "  +1.66%  "
Set-TextValue "D19" "244.16"
$ws.Range("E19").Value = "  +0.22%  "
Set-TextValue "D20" "0.000007720"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("E21").Value = "  -0.05%  "
Set-TextValue "D22" "0.9999"
$ws.Range("E22").Value = "  -0.16%  "
Set-TextValue "D23" "6.958"
$ws.Range("E23").Value = "  +4.68%  "
Set-TextValue "D24" "9.250"
$ws.Range("E24").Value = "  -1.72%  "
Set-TextValue "D25" "165.32"
$ws.Range("E25").Value = "  -0.23%  "
Set-TextValue "D26" "18.75"
$ws.Range("E26").Value = "  -1.16%  "
Set-TextValue "D27" "0.1287"
$ws.Range("E27").Value = "  +1.36%  "
Set-TextValue "D28" "2.034"
$ws.Range("E28").Value = "  -3.38%  "
$ws.Range("E29").Value = "  +3.17%  "
Set-TextValue "D30" "1.515"
$ws.Range("E30").Value = "  -2.03%  "
Set-TextValue "D31" "4.277"
$ws.Range("E31").Value = "  -0.95%  "
Set-TextValue "D32" "4.025"
$ws.Range("E32").Value = "  -1.29%  "
Set-TextValue "D33" "0.05334"
$ws.Range("E33").Value = "  +3.04%  "
Set-TextValue "D34" "1.254"
$ws.Range("E34").Value = "  -2.97%  "
Set-TextValue "D35" "0.7247"
$ws.Range("E35").Value = "  -2.56%  "
$ws.Range("E36").Value = "  -1.70%  "
Set-TextValue "D37" "0.01908"
$ws.Range("E37").Value = "  -1.58%  "
Set-TextValue "D38" "2.786"
$ws.Range("E38").Value = "  -0.30%  "
Set-TextValue "D39" "6.176"
$ws.Range("E39").Value = "  -2.93%  "
Set-TextValue "D40" "0.4405"
$ws.Range("E40").Value = "  -1.44%  "
Set-TextValue "D41" "72.31"
$ws.Range("E41").Value = "  -3.74%  "
Set-TextValue "D42" "0.9996"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  -1.66%  "
Set-TextValue "D44" "0.8242"
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("E45").Value = "  -0.42%  "
Set-TextValue "D46" "7.552"
$ws.Range("E46").Value = "  -1.11%  "
Set-TextValue "D47" "9.774"
$ws.Range("E47").Value = "  -0.65%  "
Set-TextValue "D48" "2.045.46"
$ws.Range("E48").Value = "  -1.13%  "
Set-TextValue "D49" "36.30"
$ws.Range("E49").Value = "  -2.99%  "
$ws.Range("E50").Value = "  -0.63%  "
Set-TextValue "D51" "1.468"
$ws.Range("E51").Value = "  +0.75%  "
